# Updated cryptos list on Thu Oct 12 22:30:15 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column holds text values (e.g. "26.820.25", "0.0854") - keep the
# whole column formatted as Text so Excel doesn't silently coerce
# numeric-looking strings (leading/trailing zeros, thousand-dot strings) to
# floating point numbers when we assign them below.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2 - Bitcoin
$ws.Range("D2").Value = "26.820.25"
$ws.Range("E2").Value = "  +0.11%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.541.71"
$ws.Range("E3").Value = "  -1.77%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.23%  "

# Row 5 - BNB
$ws.Range("D5").Value = "206.13"
$ws.Range("E5").Value = "  -0.23%  "

# Row 6 - XRP
$ws.Range("E6").Value = "  -0.80%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.21%  "

# Row 8 - Cardano
$ws.Range("E8").Value = "  -0.37%  "

# Row 9 - Solana
$ws.Range("D9").Value = "21.32"
$ws.Range("E9").Value = "  -3.04%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  -0.51%  "

# Row 11 - TRON
$ws.Range("D11").Value = "0.0854"
$ws.Range("E11").Value = "  -1.04%  "

# Row 12 - WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "1.760.89"
$ws.Range("E12").Value = "  -1.69%  "

# Row 13 - WrappedEther
$ws.Range("D13").Value = "1.541.22"
$ws.Range("E13").Value = "  -1.68%  "

# Row 14 - Polkadot
$ws.Range("E14").Value = "  -1.50%  "

# Row 15 - Polygon
$ws.Range("D15").Value = "0.510"
$ws.Range("E15").Value = "  -1.13%  "

# Row 16 - WrappedBTC
$ws.Range("D16").Value = "26.815.56"
$ws.Range("E16").Value = "  +0.08%  "

# Row 17 - Litecoin
$ws.Range("D17").Value = "61.22"
$ws.Range("E17").Value = "  -0.36%  "

# Row 18 - BitcoinCash
$ws.Range("D18").Value = "214.35"
$ws.Range("E18").Value = "  -0.51%  "

# Row 19 - Chainlink
$ws.Range("E19").Value = "  -2.35%  "

# Row 20 - ShibaInu
$ws.Range("E20").Value = "  +0.77%  "

# Row 21 - Dai
$ws.Range("E21").Value = "  +0.17%  "

# Row 22 - Uniswap
$ws.Range("E22").Value = "  -2.92%  "

# Row 23 - Avalanche
$ws.Range("D23").Value = "9.14"
$ws.Range("E23").Value = "  -1.94%  "

# Row 24 - Toncoin
$ws.Range("E24").Value = "  -3.34%  "

# Row 25 - Monero
$ws.Range("D25").Value = "152.16"
$ws.Range("E25").Value = "  -0.34%  "

# Row 26 - Cosmos
$ws.Range("E26").Value = "  -2.25%  "

# Row 27 - EthereumClassic
$ws.Range("D27").Value = "14.82"
$ws.Range("E27").Value = "  -1.00%  "

# Row 28 - BinanceUSD
$ws.Range("E28").Value = "  +0.19%  "

# Row 29 - Stellar
$ws.Range("E29").Value = "  -0.95%  "

# Row 30 - was Hedera, now PancakeSwap
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").Value = "1.10"
$ws.Range("E30").Value = "  -1.15%  "

# Row 31 - was PancakeSwap, now Hedera
$ws.Range("B31").Value = "Hedera"
$ws.Range("C31").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D31").Value = "0.0459"
$ws.Range("E31").Value = "  -1.56%  "

# Row 32 - Filecoin
$ws.Range("E32").Value = "  +1.91%  "

# Row 33 - Maker
$ws.Range("D33").Value = "1.368.33"
$ws.Range("E33").Value = "  -2.00%  "

# Row 34 - InternetComputer(DFINITY)
$ws.Range("E34").Value = "  +0.02%  "

# Row 35 - LidoDAOToken
$ws.Range("D35").Value = "1.50"
$ws.Range("E35").Value = "  -1.76%  "

# Row 36 - TrustWalletToken
$ws.Range("D36").Value = "0.965"
$ws.Range("E36").Value = "  +3.57%  "

# Row 37 - HuobiToken
$ws.Range("E37").Value = "  +0.10%  "

# Row 38 - VeChain
$ws.Range("E38").Value = "  +0.82%  "

# Row 39 - ImmutableX
$ws.Range("D39").Value = "0.520"
$ws.Range("E39").Value = "  -1.48%  "

# Row 40 - ARBITRUM
$ws.Range("E40").Value = "  -1.58%  "

# Row 41 - FraxShare
$ws.Range("E41").Value = "  +8.06%  "

# Row 42 - WEMIXToken
$ws.Range("D42").Value = "0.991"
$ws.Range("E42").Value = "  +0.26%  "

# Row 43 - MXToken
$ws.Range("E43").Value = "  +1.05%  "

# Row 44 - Aave
$ws.Range("D44").Value = "62.99"
$ws.Range("E44").Value = "  -0.68%  "

# Row 45 - RenderToken
$ws.Range("E45").Value = "  -3.76%  "

# Row 46 - RocketPoolETH
$ws.Range("D46").Value = "1.675.40"
$ws.Range("E46").Value = "  -1.63%  "

# Row 47 - Quant
$ws.Range("D47").Value = "84.25"
$ws.Range("E47").Value = "  -1.95%  "

# Row 48 - Cronos
$ws.Range("E48").Value = "  +3.45%  "

# Row 49 - BabyDogeCoin
$ws.Range("D49").Value = "0.0₇0980"
$ws.Range("E49").Value = "  -0.46%  "

# Row 50 - was Algorand, now USDD
$ws.Range("B50").Value = "USDD"
$ws.Range("C50").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Range("D50").Value = "1.00"
$ws.Range("E50").Value = "  +0.36%  "

# Row 51 - was USDD, now Algorand
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").Value = "0.0941"
$ws.Range("E51").Value = "  -1.41%  "
